$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.026.25"
$ws.Range("E2").Value = "'  +7.95%  "

$ws.Range("D3").Value = "'1.823.15"
$ws.Range("E3").Value = "'  +5.38%  "

$ws.Range("D4").Value = "'0.9996"

$ws.Range("D5").Value = "'246.34"
$ws.Range("E5").Value = "'  +2.41%  "

$ws.Range("D6").Value = "'0.9997"
$ws.Range("E6").Value = "'  +0.01%  "

$ws.Range("D7").Value = "'0.4924"
$ws.Range("E7").Value = "'  +1.94%  "

$ws.Range("B8").Value = "'Cardano"
$ws.Range("C8").Value = "'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").Value = "'0.2767"
$ws.Range("E8").Value = "'  +6.46%  "

$ws.Range("B9").Value = "'Dogecoin"
$ws.Range("C9").Value = "'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "'0.06379"
$ws.Range("E9").Value = "'  +3.34%  "

$ws.Range("B10").Value = "'WrappedEther"
$ws.Range("C10").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D10").Value = "'1.818.49"
$ws.Range("E10").Value = "'  +5.07%  "

$ws.Range("B11").Value = "'Solana"
$ws.Range("C11").Value = "'https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D11").Value = "'16.60"
$ws.Range("E11").Value = "'  +3.53%  "

$ws.Range("B12").Value = "'TRON"
$ws.Range("C12").Value = "'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "'0.07065"
$ws.Range("E12").Value = "'  +2.70%  "

$ws.Range("B13").Value = "'Polygon"
$ws.Range("C13").Value = "'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").Value = "'0.6428"
$ws.Range("E13").Value = "'  +6.60%  "

$ws.Range("B14").Value = "'Litecoin"
$ws.Range("C14").Value = "'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D14").Value = "'84.03"
$ws.Range("E14").Value = "'  +9.15%  "

$ws.Range("B15").Value = "'Polkadot"
$ws.Range("C15").Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "'4.692"
$ws.Range("E15").Value = "'  +5.15%  "

$ws.Range("B16").Value = "'WrappedBTC"
$ws.Range("C16").Value = "'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "'29.035.45"
$ws.Range("E16").Value = "'  +8.86%  "

$ws.Range("B17").Value = "'Dai"
$ws.Range("C17").Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D17").Value = "'0.9993"
$ws.Range("E17").Value = "'  -0.06%  "

$ws.Range("B18").Value = "'ShibaInu"
$ws.Range("C18").Value = "'https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "'0.000007300"
$ws.Range("E18").Value = "'  +2.60%  "

$ws.Range("B19").Value = "'BinanceUSD"
$ws.Range("C19").Value = "'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D19").Value = "'0.9993"
$ws.Range("E19").Value = "'  +0.02%  "

$ws.Range("B20").Value = "'Avalanche"
$ws.Range("C20").Value = "'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").Value = "'12.18"
$ws.Range("E20").Value = "'  +7.19%  "

$ws.Range("B21").Value = "'WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "'2.055.15"
$ws.Range("E21").Value = "'  +5.35%  "

$ws.Range("B22").Value = "'Uniswap"
$ws.Range("C22").Value = "'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "'4.543"
$ws.Range("E22").Value = "'  +3.39%  "

$ws.Range("B23").Value = "'Cosmos"
$ws.Range("C23").Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D23").Value = "'8.830"
$ws.Range("E23").Value = "'  +4.96%  "

$ws.Range("B24").Value = "'Chainlink"
$ws.Range("C24").Value = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D24").Value = "'5.363"
$ws.Range("E24").Value = "'  +6.02%  "

$ws.Range("B25").Value = "'Monero"
$ws.Range("C25").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").Value = "'143.41"
$ws.Range("E25").Value = "'  +2.65%  "

$ws.Range("B26").Value = "'BitcoinCash"
$ws.Range("C26").Value = "'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D26").Value = "'130.00"
$ws.Range("E26").Value = "'  +21.89%  "

$ws.Range("B27").Value = "'EthereumClassic"
$ws.Range("C27").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'16.34"
$ws.Range("E27").Value = "'  +7.50%  "

$ws.Range("B28").Value = "'LidoDAOToken"
$ws.Range("C28").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "'1.880"
$ws.Range("E28").Value = "'  +4.64%  "

$ws.Range("B29").Value = "'Toncoin"
$ws.Range("C29").Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'1.397"
$ws.Range("E29").Value = "'  +1.60%  "

$ws.Range("B30").Value = "'InternetComputer(DFINITY)"
$ws.Range("C30").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "'4.122"
$ws.Range("E30").Value = "'  +4.42%  "

$ws.Range("B31").Value = "'Stellar"
$ws.Range("C31").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").Value = "'0.08358"
$ws.Range("E31").Value = "'  +5.46%  "

$ws.Range("B32").Value = "'Filecoin"
$ws.Range("C32").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'3.765"
$ws.Range("E32").Value = "'  +2.75%  "

$ws.Range("B33").Value = "'Hedera"
$ws.Range("C33").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.04946"
$ws.Range("E33").Value = "'  +7.67%  "

$ws.Range("B34").Value = "'ARBITRUM"
$ws.Range("C34").Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Value = "'1.094"
$ws.Range("E34").Value = "'  +9.42%  "

$ws.Range("B35").Value = "'HuobiToken"
$ws.Range("C35").Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "'2.702"
$ws.Range("E35").Value = "'  +4.23%  "

$ws.Range("B36").Value = "'ImmutableX"
$ws.Range("C36").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'0.6689"
$ws.Range("E36").Value = "'  +8.50%  "

$ws.Range("B37").Value = "'RenderToken"
$ws.Range("C37").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").Value = "'2.287"
$ws.Range("E37").Value = "'  +15.17%  "

$ws.Range("B38").Value = "'MXToken"
$ws.Range("C38").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").Value = "'2.682"
$ws.Range("E38").Value = "'  +8.60%  "

$ws.Range("B39").Value = "'TrustWalletToken"
$ws.Range("C39").Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "'0.9476"
$ws.Range("E39").Value = "'  +2.50%  "

$ws.Range("B40").Value = "'FraxShare"
$ws.Range("C40").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'6.151"
$ws.Range("E40").Value = "'  +7.57%  "

$ws.Range("B41").Value = "'VeChain"
$ws.Range("C41").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "'0.01582"
$ws.Range("E41").Value = "'  +5.63%  "

$ws.Range("B42").Value = "'PaxDollar"
$ws.Range("C42").Value = "'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").Value = "'0.9996"
$ws.Range("E42").Value = "'  +0.06%  "

$ws.Range("B43").Value = "'Quant"
$ws.Range("C43").Value = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "'101.38"
$ws.Range("E43").Value = "'  +1.45%  "

$ws.Range("B44").Value = "'TheSandbox"
$ws.Range("C44").Value = "'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "'0.4060"
$ws.Range("E44").Value = "'  +5.98%  "

$ws.Range("B45").Value = "'Aptos"
$ws.Range("C45").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").Value = "'7.178"
$ws.Range("E45").Value = "'  +6.16%  "

$ws.Range("B46").Value = "'Algorand"
$ws.Range("C46").Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").Value = "'0.1217"
$ws.Range("E46").Value = "'  +5.63%  "

$ws.Range("B47").Value = "'Cronos"
$ws.Range("C47").Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "'0.05536"
$ws.Range("E47").Value = "'  +3.21%  "

$ws.Range("B48").Value = "'EnergySwap"
$ws.Range("C48").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'8.102"
$ws.Range("E48").Value = "'  +2.39%  "

$ws.Range("B49").Value = "'Elrond"
$ws.Range("C49").Value = "'https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "'31.66"
$ws.Range("E49").Value = "'  +5.11%  "

$ws.Range("B50").Value = "'NEARProtocol"
$ws.Range("C50").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "'1.298"
$ws.Range("E50").Value = "'  +4.63%  "

$ws.Range("B51").Value = "'Decentraland"
$ws.Range("C51").Value = "'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").Value = "'0.3586"
$ws.Range("E51").Value = "'  +7.06%  "
